# Automatic update of files.
# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2-28) on the active sheet from 45519 to 45520 (i.e. advance the
# date stored in column C by one day), matching the values already
# present in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45519) {
        $cell.Value = 45520
    }
}
